$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.900.85"
$ws.Range("E2").Value = "  +1.18%  "
$ws.Range("D3").Value = "2.352.96"
$ws.Range("E3").Value = "  +0.08%  "
$ws.Range("E4").Value = "  -0.22%  "
$ws.Range("D5").Value = "'0.666"
$ws.Range("E5").Value = "  +2.56%  "
$ws.Range("D6").Value = "'237.37"
$ws.Range("E6").Value = "  +2.28%  "
$ws.Range("D7").Value = "'73.21"
$ws.Range("E7").Value = "  +11.25%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  +19.49%  "
$ws.Range("D10").Value = "'0.0991"
$ws.Range("E10").Value = "  +3.29%  "
$ws.Range("D11").Value = "'28.50"
$ws.Range("E11").Value = "  +6.57%  "
$ws.Range("E12").Value = "  +2.43%  "
$ws.Range("D13").Value = "2.700.79"
$ws.Range("E13").Value = "  +0.34%  "
$ws.Range("D14").Value = "'16.66"
$ws.Range("E14").Value = "  +8.06%  "
$ws.Range("D15").Value = "'6.65"
$ws.Range("E15").Value = "  +6.29%  "
$ws.Range("E16").Value = "  +5.43%  "
$ws.Range("D17").Value = "2.359.63"
$ws.Range("E17").Value = "  +0.10%  "
$ws.Range("D18").Value = "43.784.85"
$ws.Range("E18").Value = "  +0.91%  "
$ws.Range("E19").Value = "  +3.32%  "
$ws.Range("D20").Value = "'77.87"
$ws.Range("E20").Value = "  +5.08%  "
$ws.Range("D21").Value = "'6.42"
$ws.Range("E21").Value = "  +3.25%  "
$ws.Range("D22").Value = "'253.95"
$ws.Range("E22").Value = "  +1.88%  "
$ws.Range("E23").Value = "  -0.06%  "
$ws.Range("E24").Value = "  -2.76%  "
$ws.Range("E25").Value = "  +3.13%  "
$ws.Range("D26").Value = "'10.55"
$ws.Range("E26").Value = "  +6.35%  "
$ws.Range("E27").Value = "  +0.07%  "
$ws.Range("D28").Value = "'22.40"
$ws.Range("E28").Value = "  +0.89%  "
$ws.Range("D29").Value = "'172.68"
$ws.Range("E29").Value = "  -1.32%  "
$ws.Range("E30").Value = "  +5.28%  "
$ws.Range("E31").Value = "  +1.92%  "
$ws.Range("D32").Value = "'0.132"
$ws.Range("E32").Value = "  +4.93%  "
$ws.Range("D33").Value = "'5.17"
$ws.Range("E33").Value = "  +3.66%  "
$ws.Range("D34").Value = "'0.0713"
$ws.Range("E34").Value = "  +3.67%  "
$ws.Range("D35").Value = "'5.20"
$ws.Range("E35").Value = "  +4.51%  "
$ws.Range("D36").Value = "'4.07"
$ws.Range("E36").Value = "  +12.45%  "
$ws.Range("E37").Value = "  -4.56%  "
$ws.Range("E38").Value = "  -1.08%  "
$ws.Range("E39").Value = "  +6.32%  "
$ws.Range("D40").Value = "'19.64"
$ws.Range("E40").Value = "  +8.71%  "
$ws.Range("E41").Value = "  -0.20%  "
$ws.Range("D42").Value = "'8.82"
$ws.Range("E42").Value = "  -2.29%  "
$ws.Range("E43").Value = "  +3.91%  "
$ws.Range("D44").Value = "'0.0981"
$ws.Range("E44").Value = "  +3.42%  "
$ws.Range("E45").Value = "  -1.47%  "
$ws.Range("D46").Value = "'4.45"
$ws.Range("E46").Value = "  +1.27%  "
$ws.Range("D47").Value = "'97.86"
$ws.Range("E48").Value = "  +11.13%  "
$ws.Range("D49").Value = "'2.32"
$ws.Range("E49").Value = "  +2.84%  "
$ws.Range("D50").Value = "1.435.92"
$ws.Range("E50").Value = "  -0.30%  "
$ws.Range("E51").Value = "  +1.49%  "
